# Deploying to gh-pages: add 2021 / 2022 (columns N, O) data to the
# "Percentage of population with sustainable access to sanitation
# facilities" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (thin separator row under the year headers): extend the
#     existing bottom-border formatting (as used by K3:M3) into N3:O3.
[void]$ws.Range("K3").Copy()
[void]$ws.Range("N3:O3").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 4 (year headers 2011..2020 -> add 2021, 2022): copy M4's
#     format (year-header style) into N4:O4, then set the year values.
[void]$ws.Range("M4").Copy()
[void]$ws.Range("N4:O4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N4").Value = 2021
$ws.Range("O4").Value = 2022

# --- Row 5 (Kyrgyz Republic total, bold): numeric "0.0" format, bold.
#     Style N5 directly (this mutates its auto-assigned xf in place, with
#     no leftover intermediate style since it is the sole owner), then
#     clone that exact cell format onto O5 via a formats-only paste so
#     both cells end up sharing one new style record.
$ws.Range("N5").Font.Bold = $true
$ws.Range("N5").NumberFormat = "0.0"
[void]$ws.Range("N5").Copy()
[void]$ws.Range("O5").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N5").Value = 40.007977647471066
$ws.Range("O5").Value = 42.620582506455563

# --- Rows 6-13 (oblast detail rows): numeric "0.0" format, regular weight.
$ws.Range("N6:O13").NumberFormat = "0.0"

$ws.Range("N6").Value = 5.7072514621689896
$ws.Range("O6").Value = 8.1443914479075037

$ws.Range("N7").Value = 8.9893229854028949
$ws.Range("O7").Value = 10.715961386284755

$ws.Range("N8").Value = 66.307512472824584
$ws.Range("O8").Value = 81.977461999426666

$ws.Range("N9").Value = 23.475213049310256
$ws.Range("O9").Value = 29.828871240443185

$ws.Range("N10").Value = 9.8045372040896162
$ws.Range("O10").Value = 9.7218425128664112

$ws.Range("N11").Value = 9.3737779268960448
$ws.Range("O11").Value = 8.6167819403064012

$ws.Range("N12").Value = 70.457032471318783
$ws.Range("O12").Value = 69.915337594090886

$ws.Range("N13").Value = 98.411252120183207
$ws.Range("O13").Value = 99.08571752721997

# --- Row 14 (Osh city, bottom-bordered): copy K3's bottom border into
#     N14:O14, then apply the "0.0" number format and values.
[void]$ws.Range("K3").Copy()
[void]$ws.Range("N14:O14").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("N14:O14").NumberFormat = "0.0"
$ws.Range("N14").Value = 63.900563564170795
$ws.Range("O14").Value = 64.805252627098838

# --- Move the active selection to P8 (matches the recorded cursor
#     position after the edit).
[void]$ws.Range("P8").Select()
